# Task 7 at Project_Playground
# Fill in the newly tracked task row (row 7) on the PM-SHEET with the
# task "Informing abou Android Std." assigned to "Tumfart/Trimbacher".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PM-SHEET")

# Responsible person first, then task name, so the shared-string table
# is populated in the same order as the authored workbook.
$ws.Range("G7").Value = "Tumfart/Trimbacher"
$ws.Range("B7").Value = "Informing abou Android Std."
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = 20
$ws.Range("E7").Value = 14
$ws.Range("F7").Formula = "=D7-E7"

# Leave the cursor where the author left it when saving.
$ws.Range("E7").Select()
